$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/15/2024  Through  7/21/2024"

# --- Cells that change from a placeholder text style to a numeric style ---
$ws.Range("M14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M14").Value = -100
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 2
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = -100
$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("C16").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 2
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = -100

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 7
$ws.Range("K15").Value = -28.571428571428
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 6
$ws.Range("H16").Value = -45.454545454545
$ws.Range("I16").Value = 70
$ws.Range("J16").Value = 58
$ws.Range("K16").Value = 20.689655172413
$ws.Range("L16").Value = 14.754098360655
$ws.Range("M16").Value = -46.153846153846
$ws.Range("N16").Value = -81.958762886597
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 90
$ws.Range("J17").Value = 97
$ws.Range("K17").Value = -7.216494845360
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 34.328358208955
$ws.Range("N17").Value = -49.438202247191
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("I18").Value = 36
$ws.Range("J18").Value = 52
$ws.Range("K18").Value = -30.769230769230
$ws.Range("L18").Value = -37.931034482758
$ws.Range("M18").Value = -73.722627737226
$ws.Range("N18").Value = -94.642857142857
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -60
$ws.Range("G19").Value = 65
$ws.Range("H19").Value = -47.692307692307
$ws.Range("I19").Value = 313
$ws.Range("J19").Value = 350
$ws.Range("K19").Value = -10.571428571428
$ws.Range("L19").Value = 15.925925925925
$ws.Range("M19").Value = 23.228346456692
$ws.Range("N19").Value = -14.010989010989
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = -35
$ws.Range("I20").Value = 87
$ws.Range("J20").Value = 69
$ws.Range("K20").Value = 26.086956521739
$ws.Range("L20").Value = 35.9375
$ws.Range("M20").Value = -5.434782608695
$ws.Range("N20").Value = -94.419499679281
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -37.931034482758
$ws.Range("F21").Value = 75
$ws.Range("G21").Value = 121
$ws.Range("H21").Value = -38.016528925619
$ws.Range("I21").Value = 601
$ws.Range("J21").Value = 635
$ws.Range("K21").Value = -5.354330708661
$ws.Range("L21").Value = 9.272727272727
$ws.Range("M21").Value = -12.772133526850
$ws.Range("N21").Value = -81.148055207026
$ws.Range("D23").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("H23").Value = -80
$ws.Range("J23").Value = 14
$ws.Range("K23").Value = 42.857142857142
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -41.935483870967
$ws.Range("F24").Value = 122
$ws.Range("G24").Value = 101
$ws.Range("H24").Value = 20.792079207920
$ws.Range("I24").Value = 892
$ws.Range("J24").Value = 675
$ws.Range("K24").Value = 32.148148148148
$ws.Range("L24").Value = 45.990180032733
$ws.Range("M24").Value = 77.689243027888
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = -39.130434782608
$ws.Range("F25").Value = 102
$ws.Range("G25").Value = 79
$ws.Range("H25").Value = 29.113924050632
$ws.Range("I25").Value = 735
$ws.Range("J25").Value = 468
$ws.Range("K25").Value = 57.051282051282
$ws.Range("L25").Value = 93.931398416886
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 25
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = 4.166666666666
$ws.Range("I26").Value = 180
$ws.Range("J26").Value = 148
$ws.Range("K26").Value = 21.621621621621
$ws.Range("L26").Value = 44
$ws.Range("M26").Value = -7.216494845360
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -66.666666666666
$ws.Range("J27").Value = 15
$ws.Range("K27").Value = -53.333333333333
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 20
$ws.Range("J28").Value = 18
$ws.Range("K28").Value = 11.111111111111
$ws.Range("L28").Value = 5.263157894736
$ws.Range("M29").Value = -75
$ws.Range("M30").Value = -71.428571428571
